# Apply the "Fix: Portfolio fetching (Batch MultiIndex + NaN handling) & Firebase Auth"
# update to the financial/physics roadmap tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row (row 1): strip the bold / bordered / centered style that was
#    applied to A1:N1 so the header goes back to the plain default style.
# ---------------------------------------------------------------------------
$ws.Range("A1:N1").ClearFormats()

# ---------------------------------------------------------------------------
# 2) Add five new tracker columns (O:S) with their header labels.
# ---------------------------------------------------------------------------
$ws.Range("O1").Value = "Area"
$ws.Range("P1").Value = "Task"
$ws.Range("Q1").Value = "Priority"
$ws.Range("R1").Value = "Status"
$ws.Range("S1").Value = "Notes"

# ---------------------------------------------------------------------------
# 3) Rewrite row 7 (previously the duplicated "Fourier UI Control & Window"
#    entry, ID 6) to describe the portfolio loading / Firebase auth fix.
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "Fix Portfolio Loading Freeze & Price Data NaNs"
$ws.Range("C7").Value = "Backend / Infrastructure"
$ws.Range("D7:F7").ClearContents()
$ws.Range("G7").Value = "Done"
$ws.Range("H7").Value = "Critical"
$ws.Range("I7:M7").ClearContents()
$ws.Range("N7").Value = "Resolved Firebase Auth (JWT) error, fixed yfinance batch fetching for MultiIndex, handled NaN prices for US stocks. 2026-01-30"

# ---------------------------------------------------------------------------
# 4) Remove the old row 8 ("Media Portafoglio Equipesata", ID 7) entirely -
#    the tracker now ends at row 7.
# ---------------------------------------------------------------------------
$ws.Rows("8").Delete()
